$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Update "Objetivos:" (row 10) body text to the new Portuguese text
# ------------------------------------------------------------------
$ws.Range("B10").Value = "Complementar a formação dos estudantes abordando, com maior profundidade, tópicos atuais e relevantes e atualizar com temas no estado da arte."
$ws.Range("C10").Value = "Complementar a formação dos estudantes abordando, com maior profundidade, tópicos atuais e relevantes e atualizar com temas no estado da arte."

# ------------------------------------------------------------------
# 2) Insert a new blank row at 13 (pushes the existing rows 13-21,
#    together with their heights/labels, down to rows 14-22), then
#    populate the freshly inserted row 13 with the professor's name
#    in columns B/C only (column A stays empty, matching the
#    "Docentes responsáveis:" block layout).
# ------------------------------------------------------------------
$ws.Rows.Item(13).Insert()
$ws.Range("A13").Clear()
$ws.Range("B13").Value = "5817650 - Érica Leonor Romão"
$ws.Range("C13").Value = "5817650 - Érica Leonor Romão"
$ws.Range("B2").Copy()
$ws.Range("B13").PasteSpecial(-4122)
$ws.Range("C2").Copy()
$ws.Range("C13").PasteSpecial(-4122)

# ------------------------------------------------------------------
# 3) Update remaining body texts that changed content (labels already
#    shifted correctly by the row insert above).
# ------------------------------------------------------------------
# "Programa resumido:" (now row 14)
$ws.Range("B14").Value = "A definir, de acordo com o tópico programado"
$ws.Range("C14").Value = "A definir, de acordo com o tópico programado"

# "Programa:" (now row 16)
$ws.Range("B16").Value = "O conteúdo desta disciplina optativa será de acordo com o tópico a ser programado, devendo abordar assuntos complementares ao conteúdo regular do curso de graduação."
$ws.Range("C16").Value = "O conteúdo desta disciplina optativa será de acordo com o tópico a ser programado, devendo abordar assuntos complementares ao conteúdo regular do curso de graduação."

# "Método:" (now row 19)
$ws.Range("B19").Value = "Esta disciplina deverá conter no mínimo duas avaliações denominadas P1 e P2. A P2 deverá englobar toda a matéria ministrada ao longo do semestre, abrangendo todos os tópicos previstos na ementa. As avalições podem ser: seminários, trabalhos, projetos ou outra forma de avaliação definida pelo professor. Sendo necessário no mínimo uma avaliação na forma de prova escrita."
$ws.Range("C19").Value = "Esta disciplina deverá conter no mínimo duas avaliações denominadas P1 e P2. A P2 deverá englobar toda a matéria ministrada ao longo do semestre, abrangendo todos os tópicos previstos na ementa. As avalições podem ser: seminários, trabalhos, projetos ou outra forma de avaliação definida pelo professor. Sendo necessário no mínimo uma avaliação na forma de prova escrita."

# "Critério:" (now row 20)
$ws.Range("B20").Value = "Média ponderada das avaliações (M)."
$ws.Range("C20").Value = "Média ponderada das avaliações (M)."

# "Norma de recuperação:" (now row 21)
$ws.Range("B21").Value = "A recuperação será composta por uma única prova (RC) englobando toda a matéria ministrada ao longo do semestre. A média final, para os alunos em recuperação, será calculada com base na relação: MF=(M+RC)/2"
$ws.Range("C21").Value = "A recuperação será composta por uma única prova (RC) englobando toda a matéria ministrada ao longo do semestre. A média final, para os alunos em recuperação, será calculada com base na relação: MF=(M+RC)/2"

# ------------------------------------------------------------------
# 4) Add the brand-new "Bibliografia:" content row (row 22), copying
#    formatting/row-height from the row above (row 21).
# ------------------------------------------------------------------
$ws.Range("B22").Value = "Livros, artigos ou texto fornecido pelo docente responsável extraídos de livros ou revistas especializadas na área de Meio Ambiente."
$ws.Range("C22").Value = "Livros, artigos ou texto fornecido pelo docente responsável extraídos de livros ou revistas especializadas na área de Meio Ambiente."
$ws.Range("B21").Copy()
$ws.Range("B22").PasteSpecial(-4122)
$ws.Range("C21").Copy()
$ws.Range("C22").PasteSpecial(-4122)
$ws.Rows.Item(22).RowHeight = 120
